$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: apply the "header" style already used throughout this workbook
# (bold font, thin border all around, centered horizontally, top-aligned
# vertically -- cellXfs index 1 in the original file) to a destination
# Range by copying the formatting from a cell that already carries it
# (the "entity" header cell on the "industry" sheet, which this edit never
# removes). Using copy/paste-special for formats (instead of touching
# Font/Borders/Alignment individually) avoids creating any new style
# definitions, keeping styles.xml unchanged.
# ---------------------------------------------------------------------------
$styleSource = $wb.Worksheets.Item("industry").Range("A1")

function Set-HeaderStyle($range) {
    $styleSource.Copy()
    $range.PasteSpecial(-4122)   # xlPasteFormats
}

# ===========================================================================
# Sheet "industry": drop the two rows for entities 1 and 2 (FoodWasteBuyer1/2)
# leaving only entities 262 and 261 (FoodWasteSeller2/1).
# ===========================================================================
$wsIndustry = $wb.Worksheets.Item("industry")
$wsIndustry.Rows.Item(2).Delete()
$wsIndustry.Rows.Item(2).Delete()

# ===========================================================================
# Sheet "entity": replace rows 2-5 content; entities become 262, 261, then
# the brand-new entities E6, E7, E8 (one new row added).
# ===========================================================================
$wsEntity = $wb.Worksheets.Item("entity")
$wsEntity.Range("A2").Value = 262
$wsEntity.Range("A3").Value = 261
$wsEntity.Range("A4").Value = "E6"
$wsEntity.Range("A5").Value = "E7"
$wsEntity.Range("A6").Value = "E8"

# ===========================================================================
# Sheet "industry_demand": remove the two data rows, keep only the header.
# ===========================================================================
$wsIndustryDemand = $wb.Worksheets.Item("industry_demand")
$wsIndustryDemand.Rows.Item(2).Delete()
$wsIndustryDemand.Rows.Item(2).Delete()

# ===========================================================================
# Rename the old "feasible" sheet to "investment_demand" and give it new
# content: entity / material / quantity / reserve_price for E6, E7, E8.
# ===========================================================================
$wsInvestmentDemand = $wb.Worksheets.Item("feasible")
$wsInvestmentDemand.Name = "investment_demand"
$wsInvestmentDemand.Cells.Clear()

$wsInvestmentDemand.Range("A1").Value = "entity"
$wsInvestmentDemand.Range("B1").Value = "material"
$wsInvestmentDemand.Range("C1").Value = "quantity"
$wsInvestmentDemand.Range("D1").Value = "reserve_price"
Set-HeaderStyle($wsInvestmentDemand.Range("A1:D1"))

$wsInvestmentDemand.Range("A2").Value = "E6"
$wsInvestmentDemand.Range("B2").Value = "Food Waste"
$wsInvestmentDemand.Range("C2").Value = 100
$wsInvestmentDemand.Range("D2").Value = 100

$wsInvestmentDemand.Range("A3").Value = "E7"
$wsInvestmentDemand.Range("B3").Value = "Food Waste"
$wsInvestmentDemand.Range("C3").Value = 100
$wsInvestmentDemand.Range("D3").Value = 90

$wsInvestmentDemand.Range("A4").Value = "E8"
$wsInvestmentDemand.Range("B4").Value = "Food Waste"
$wsInvestmentDemand.Range("C4").Value = 100
$wsInvestmentDemand.Range("D4").Value = 100

# ===========================================================================
# Brand-new sheet "invest_cost" right after "investment_demand": entity /
# invest_cost for E6, E7, E8.
# ===========================================================================
$wsInvestCost = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsInvestmentDemand)
$wsInvestCost.Name = "invest_cost"

$wsInvestCost.Range("A1").Value = "entity"
$wsInvestCost.Range("B1").Value = "invest_cost"
Set-HeaderStyle($wsInvestCost.Range("A1:B1"))

$wsInvestCost.Range("A2").Value = "E6"
$wsInvestCost.Range("B2").Value = 30000

$wsInvestCost.Range("A3").Value = "E7"
$wsInvestCost.Range("B3").Value = 35000

$wsInvestCost.Range("A4").Value = "E8"
$wsInvestCost.Range("B4").Value = 30000

# ===========================================================================
# Brand-new sheet "feasible" right after "invest_cost" (the old "feasible"
# sheet was renamed to "investment_demand" above, so this re-creates a
# fresh "feasible" sheet at the correct position with the reset template
# content: header row 0..4, rows for 262/261 with zeros, and placeholder
# rows for the new entities E6/E7/E8).
# ===========================================================================
$wsFeasible = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsInvestCost)
$wsFeasible.Name = "feasible"

$wsFeasible.Range("A1").Value = "entity"
$wsFeasible.Range("B1").Value = 0
$wsFeasible.Range("C1").Value = 1
$wsFeasible.Range("D1").Value = 2
$wsFeasible.Range("E1").Value = 3
$wsFeasible.Range("F1").Value = 4
Set-HeaderStyle($wsFeasible.Range("A1:F1"))

$wsFeasible.Range("A2").Value = 262
Set-HeaderStyle($wsFeasible.Range("A2"))
$wsFeasible.Range("D2").Value = 0
$wsFeasible.Range("E2").Value = 0
$wsFeasible.Range("F2").Value = 0

$wsFeasible.Range("A3").Value = 261
Set-HeaderStyle($wsFeasible.Range("A3"))
$wsFeasible.Range("D3").Value = 0
$wsFeasible.Range("E3").Value = 0
$wsFeasible.Range("F3").Value = 0

$wsFeasible.Range("A4").Value = "E6"
Set-HeaderStyle($wsFeasible.Range("A4"))

$wsFeasible.Range("A5").Value = "E7"
Set-HeaderStyle($wsFeasible.Range("A5"))

$wsFeasible.Range("A6").Value = "E8"
Set-HeaderStyle($wsFeasible.Range("A6"))

# ===========================================================================
# Existing "distance" sheet: reset to the same template content as the new
# "feasible" sheet (it loses its previous distance values and now shares
# the same reset structure: header 0..4, zeros for 262/261, placeholder
# rows for E6/E7/E8).
# ===========================================================================
$wsDistance = $wb.Worksheets.Item("distance")
$wsDistance.Cells.Clear()

$wsDistance.Range("A1").Value = "entity"
$wsDistance.Range("B1").Value = 0
$wsDistance.Range("C1").Value = 1
$wsDistance.Range("D1").Value = 2
$wsDistance.Range("E1").Value = 3
$wsDistance.Range("F1").Value = 4
Set-HeaderStyle($wsDistance.Range("A1:F1"))

$wsDistance.Range("A2").Value = 262
Set-HeaderStyle($wsDistance.Range("A2"))
$wsDistance.Range("D2").Value = 0
$wsDistance.Range("E2").Value = 0
$wsDistance.Range("F2").Value = 0

$wsDistance.Range("A3").Value = 261
Set-HeaderStyle($wsDistance.Range("A3"))
$wsDistance.Range("D3").Value = 0
$wsDistance.Range("E3").Value = 0
$wsDistance.Range("F3").Value = 0

$wsDistance.Range("A4").Value = "E6"
Set-HeaderStyle($wsDistance.Range("A4"))

$wsDistance.Range("A5").Value = "E7"
Set-HeaderStyle($wsDistance.Range("A5"))

$wsDistance.Range("A6").Value = "E8"
Set-HeaderStyle($wsDistance.Range("A6"))

$excel.CutCopyMode = $false

Write-Output "edit complete"
